$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15 for LA (shifts MD..WY down by one row)
$ws.Rows.Item(15).Insert()

# Refresh all data rows (2-40) with the updated 12-days-before-election figures
# Row 2: AK
$ws.Cells.Item(2,1).Value = "AK"
$ws.Cells.Item(2,2).Value = 40537
$ws.Cells.Item(2,3).Value = 24307
$ws.Cells.Item(2,4).Value = 126522
$ws.Cells.Item(2,5).Value = 45889
$ws.Cells.Item(2,6).Value = 85985
$ws.Cells.Item(2,7).Value = 212.11
$ws.Cells.Item(2,8).Value = 21582
$ws.Cells.Item(2,9).Value = 88.79000000000001

# Row 3: AR
$ws.Cells.Item(3,1).Value = "AR"
$ws.Cells.Item(3,2).Value = 215995
$ws.Cells.Item(3,3).Value = 215980
$ws.Cells.Item(3,4).Value = 114063
$ws.Cells.Item(3,5).Value = 68779
$ws.Cells.Item(3,6).Value = -101932
$ws.Cells.Item(3,7).Value = -47.19
$ws.Cells.Item(3,8).Value = -147201
$ws.Cells.Item(3,9).Value = -68.15000000000001

# Row 4: AZ
$ws.Cells.Item(4,1).Value = "AZ"
$ws.Cells.Item(4,2).Value = 2332640
$ws.Cells.Item(4,3).Value = 919654
$ws.Cells.Item(4,4).Value = 3233099
$ws.Cells.Item(4,5).Value = 927523
$ws.Cells.Item(4,6).Value = 900459
$ws.Cells.Item(4,7).Value = 38.6
$ws.Cells.Item(4,8).Value = 7869
$ws.Cells.Item(4,9).Value = 0.86

# Row 5: CA
$ws.Cells.Item(5,1).Value = "CA"
$ws.Cells.Item(5,2).Value = 10964487
$ws.Cells.Item(5,3).Value = 2050993
$ws.Cells.Item(5,4).Value = 3923785
$ws.Cells.Item(5,5).Value = 3923785
$ws.Cells.Item(5,6).Value = -7040702
$ws.Cells.Item(5,7).Value = -64.20999999999999
$ws.Cells.Item(5,8).Value = 1872792
$ws.Cells.Item(5,9).Value = 91.31

# Row 6: CO
$ws.Cells.Item(6,1).Value = "CO"
$ws.Cells.Item(6,2).Value = 3166218
$ws.Cells.Item(6,3).Value = 700234
$ws.Cells.Item(6,4).Value = 917942
$ws.Cells.Item(6,5).Value = 917942
$ws.Cells.Item(6,6).Value = -2248276
$ws.Cells.Item(6,7).Value = -71.01000000000001
$ws.Cells.Item(6,8).Value = 217708
$ws.Cells.Item(6,9).Value = 31.09

# Row 7: CT
$ws.Cells.Item(7,1).Value = "CT"
$ws.Cells.Item(7,4).Value = 611171
$ws.Cells.Item(7,5).Value = 365585

# Row 8: DE
$ws.Cells.Item(8,1).Value = "DE"
$ws.Cells.Item(8,2).Value = 20538
$ws.Cells.Item(8,3).Value = 14839
$ws.Cells.Item(8,4).Value = 170167
$ws.Cells.Item(8,5).Value = 90520
$ws.Cells.Item(8,6).Value = 149629
$ws.Cells.Item(8,7).Value = 728.55
$ws.Cells.Item(8,8).Value = 75681
$ws.Cells.Item(8,9).Value = 510.01

# Row 9: FL
$ws.Cells.Item(9,1).Value = "FL"
$ws.Cells.Item(9,2).Value = 4394336
$ws.Cells.Item(9,3).Value = 2832896
$ws.Cells.Item(9,4).Value = 6066784
$ws.Cells.Item(9,5).Value = 2954163
$ws.Cells.Item(9,6).Value = 1672448
$ws.Cells.Item(9,7).Value = 38.06
$ws.Cells.Item(9,8).Value = 121267
$ws.Cells.Item(9,9).Value = 4.28

# Row 10: GA
$ws.Cells.Item(10,1).Value = "GA"
$ws.Cells.Item(10,2).Value = 1167786
$ws.Cells.Item(10,3).Value = 1066169
$ws.Cells.Item(10,4).Value = 2493877
$ws.Cells.Item(10,5).Value = 1658463
$ws.Cells.Item(10,6).Value = 1326091
$ws.Cells.Item(10,7).Value = 113.56
$ws.Cells.Item(10,8).Value = 592294
$ws.Cells.Item(10,9).Value = 55.55

# Row 11: IA
$ws.Cells.Item(11,1).Value = "IA"
$ws.Cells.Item(11,2).Value = 522200
$ws.Cells.Item(11,3).Value = 378648
$ws.Cells.Item(11,4).Value = 785173
$ws.Cells.Item(11,5).Value = 511840
$ws.Cells.Item(11,6).Value = 262973
$ws.Cells.Item(11,7).Value = 50.36
$ws.Cells.Item(11,8).Value = 133192
$ws.Cells.Item(11,9).Value = 35.18

# Row 12: ID
$ws.Cells.Item(12,1).Value = "ID"
$ws.Cells.Item(12,2).Value = 133705
$ws.Cells.Item(12,3).Value = 100585
$ws.Cells.Item(12,4).Value = 422383
$ws.Cells.Item(12,5).Value = 203525
$ws.Cells.Item(12,6).Value = 288678
$ws.Cells.Item(12,7).Value = 215.91
$ws.Cells.Item(12,8).Value = 102940
$ws.Cells.Item(12,9).Value = 102.34

# Row 13: IL
$ws.Cells.Item(13,1).Value = "IL"
$ws.Cells.Item(13,2).Value = 868029
$ws.Cells.Item(13,3).Value = 604297
$ws.Cells.Item(13,4).Value = 2343725
$ws.Cells.Item(13,5).Value = 874521
$ws.Cells.Item(13,6).Value = 1475696
$ws.Cells.Item(13,7).Value = 170.01
$ws.Cells.Item(13,8).Value = 270224
$ws.Cells.Item(13,9).Value = 44.72

# Row 14: KS
$ws.Cells.Item(14,1).Value = "KS"
$ws.Cells.Item(14,2).Value = 262582
$ws.Cells.Item(14,3).Value = 152076
$ws.Cells.Item(14,4).Value = 508299
$ws.Cells.Item(14,5).Value = 76442
$ws.Cells.Item(14,6).Value = 245717
$ws.Cells.Item(14,7).Value = 93.58
$ws.Cells.Item(14,8).Value = -75634
$ws.Cells.Item(14,9).Value = -49.73

# Row 15: LA
$ws.Cells.Item(15,1).Value = "LA"
$ws.Cells.Item(15,2).Value = 225438
$ws.Cells.Item(15,3).Value = 225410
$ws.Cells.Item(15,4).Value = 340868
$ws.Cells.Item(15,5).Value = 340868
$ws.Cells.Item(15,6).Value = 115430
$ws.Cells.Item(15,7).Value = 51.2
$ws.Cells.Item(15,8).Value = 115458
$ws.Cells.Item(15,9).Value = 51.22

# Row 16: MD
$ws.Cells.Item(16,1).Value = "MD"
$ws.Cells.Item(16,2).Value = 330760
$ws.Cells.Item(16,3).Value = 228038
$ws.Cells.Item(16,4).Value = 1597173
$ws.Cells.Item(16,5).Value = 688963
$ws.Cells.Item(16,6).Value = 1266413
$ws.Cells.Item(16,7).Value = 382.88
$ws.Cells.Item(16,8).Value = 460925
$ws.Cells.Item(16,9).Value = 202.13

# Row 17: ME
$ws.Cells.Item(17,1).Value = "ME"
$ws.Cells.Item(17,2).Value = 181635
$ws.Cells.Item(17,3).Value = 125966
$ws.Cells.Item(17,4).Value = 402352
$ws.Cells.Item(17,5).Value = 271863
$ws.Cells.Item(17,6).Value = 220717
$ws.Cells.Item(17,7).Value = 121.52
$ws.Cells.Item(17,8).Value = 145897
$ws.Cells.Item(17,9).Value = 115.82

# Row 18: MI
$ws.Cells.Item(18,1).Value = "MI"
$ws.Cells.Item(18,2).Value = 1052138
$ws.Cells.Item(18,3).Value = 617268
$ws.Cells.Item(18,4).Value = 3002444
$ws.Cells.Item(18,5).Value = 1694413
$ws.Cells.Item(18,6).Value = 1950306
$ws.Cells.Item(18,7).Value = 185.37
$ws.Cells.Item(18,8).Value = 1077145
$ws.Cells.Item(18,9).Value = 174.5

# Row 19: MN
$ws.Cells.Item(19,1).Value = "MN"
$ws.Cells.Item(19,2).Value = 252636
$ws.Cells.Item(19,3).Value = 252603
$ws.Cells.Item(19,4).Value = 1023646
$ws.Cells.Item(19,5).Value = 1023646
$ws.Cells.Item(19,6).Value = 771010
$ws.Cells.Item(19,7).Value = 305.19
$ws.Cells.Item(19,8).Value = 771043
$ws.Cells.Item(19,9).Value = 305.24

# Row 20: MT
$ws.Cells.Item(20,1).Value = "MT"
$ws.Cells.Item(20,2).Value = 323097
$ws.Cells.Item(20,3).Value = 163995
$ws.Cells.Item(20,4).Value = 630582
$ws.Cells.Item(20,5).Value = 185978
$ws.Cells.Item(20,6).Value = 307485
$ws.Cells.Item(20,7).Value = 95.17
$ws.Cells.Item(20,8).Value = 21983
$ws.Cells.Item(20,9).Value = 13.4

# Row 21: NC
$ws.Cells.Item(21,1).Value = "NC"
$ws.Cells.Item(21,2).Value = 1311769
$ws.Cells.Item(21,3).Value = 1203096
$ws.Cells.Item(21,4).Value = 1839773
$ws.Cells.Item(21,5).Value = 1829433
$ws.Cells.Item(21,6).Value = 528004
$ws.Cells.Item(21,7).Value = 40.25
$ws.Cells.Item(21,8).Value = 626337
$ws.Cells.Item(21,9).Value = 52.06

# Row 22: ND
$ws.Cells.Item(22,1).Value = "ND"
$ws.Cells.Item(22,2).Value = 77686
$ws.Cells.Item(22,3).Value = 45214
$ws.Cells.Item(22,4).Value = 194931
$ws.Cells.Item(22,5).Value = 104270
$ws.Cells.Item(22,6).Value = 117245
$ws.Cells.Item(22,7).Value = 150.92
$ws.Cells.Item(22,8).Value = 59056
$ws.Cells.Item(22,9).Value = 130.61

# Row 23: NE
$ws.Cells.Item(23,1).Value = "NE"
$ws.Cells.Item(23,2).Value = 187886
$ws.Cells.Item(23,3).Value = 109959
$ws.Cells.Item(23,4).Value = 490240
$ws.Cells.Item(23,5).Value = 285581
$ws.Cells.Item(23,6).Value = 302354
$ws.Cells.Item(23,7).Value = 160.92
$ws.Cells.Item(23,8).Value = 175622
$ws.Cells.Item(23,9).Value = 159.72

# Row 24: NJ
$ws.Cells.Item(24,1).Value = "NJ"
$ws.Cells.Item(24,2).Value = 376161
$ws.Cells.Item(24,3).Value = 214216
$ws.Cells.Item(24,4).Value = 5802187
$ws.Cells.Item(24,5).Value = 1686290
$ws.Cells.Item(24,6).Value = 5426026
$ws.Cells.Item(24,7).Value = 1442.47
$ws.Cells.Item(24,8).Value = 1472074
$ws.Cells.Item(24,9).Value = 687.1900000000001

# Row 25: NM
$ws.Cells.Item(25,1).Value = "NM"
$ws.Cells.Item(25,2).Value = 267791
$ws.Cells.Item(25,3).Value = 237733
$ws.Cells.Item(25,4).Value = 523961
$ws.Cells.Item(25,5).Value = 339043
$ws.Cells.Item(25,6).Value = 256170
$ws.Cells.Item(25,7).Value = 95.66
$ws.Cells.Item(25,8).Value = 101310
$ws.Cells.Item(25,9).Value = 42.62

# Row 26: NV
$ws.Cells.Item(26,1).Value = "NV"
$ws.Cells.Item(26,2).Value = 341914
$ws.Cells.Item(26,3).Value = 300397
$ws.Cells.Item(26,4).Value = 1479497
$ws.Cells.Item(26,5).Value = 270000
$ws.Cells.Item(26,6).Value = 1137583
$ws.Cells.Item(26,7).Value = 332.71
$ws.Cells.Item(26,8).Value = -30397
$ws.Cells.Item(26,9).Value = -10.12

# Row 27: OH
$ws.Cells.Item(27,1).Value = "OH"
$ws.Cells.Item(27,2).Value = 1473264
$ws.Cells.Item(27,3).Value = 892798
$ws.Cells.Item(27,4).Value = 2500962
$ws.Cells.Item(27,5).Value = 1264667
$ws.Cells.Item(27,6).Value = 1027698
$ws.Cells.Item(27,7).Value = 69.76000000000001
$ws.Cells.Item(27,8).Value = 371869
$ws.Cells.Item(27,9).Value = 41.65

# Row 28: OK
$ws.Cells.Item(28,1).Value = "OK"
$ws.Cells.Item(28,4).Value = 317441
$ws.Cells.Item(28,5).Value = 0

# Row 29: OR
$ws.Cells.Item(29,1).Value = "OR"
$ws.Cells.Item(29,2).Value = 444899
$ws.Cells.Item(29,3).Value = 444829
$ws.Cells.Item(29,4).Value = 274866
$ws.Cells.Item(29,5).Value = 274866
$ws.Cells.Item(29,6).Value = -170033
$ws.Cells.Item(29,7).Value = -38.22
$ws.Cells.Item(29,8).Value = -169963
$ws.Cells.Item(29,9).Value = -38.21

# Row 30: PA
$ws.Cells.Item(30,1).Value = "PA"
$ws.Cells.Item(30,4).Value = 2836018
$ws.Cells.Item(30,5).Value = 1024804

# Row 31: RI
$ws.Cells.Item(31,1).Value = "RI"
$ws.Cells.Item(31,2).Value = 33020
$ws.Cells.Item(31,3).Value = 12875
$ws.Cells.Item(31,4).Value = 216835
$ws.Cells.Item(31,5).Value = 41757
$ws.Cells.Item(31,6).Value = 183815
$ws.Cells.Item(31,7).Value = 556.68
$ws.Cells.Item(31,8).Value = 28882
$ws.Cells.Item(31,9).Value = 224.33

# Row 32: SD
$ws.Cells.Item(32,1).Value = "SD"
$ws.Cells.Item(32,4).Value = 171622
$ws.Cells.Item(32,5).Value = 137187

# Row 33: TN
$ws.Cells.Item(33,1).Value = "TN"
$ws.Cells.Item(33,2).Value = 942946
$ws.Cells.Item(33,3).Value = 942919
$ws.Cells.Item(33,4).Value = 810217
$ws.Cells.Item(33,5).Value = 810217
$ws.Cells.Item(33,6).Value = -132729
$ws.Cells.Item(33,7).Value = -14.08
$ws.Cells.Item(33,8).Value = -132702
$ws.Cells.Item(33,9).Value = -14.07

# Row 34: TX
$ws.Cells.Item(34,1).Value = "TX"
$ws.Cells.Item(34,2).Value = 1790826
$ws.Cells.Item(34,3).Value = 1762855
$ws.Cells.Item(34,4).Value = 4375786
$ws.Cells.Item(34,5).Value = 4375786
$ws.Cells.Item(34,6).Value = 2584960
$ws.Cells.Item(34,7).Value = 144.34
$ws.Cells.Item(34,8).Value = 2612931
$ws.Cells.Item(34,9).Value = 148.22

# Row 35: UT
$ws.Cells.Item(35,1).Value = "UT"
$ws.Cells.Item(35,2).Value = 1007467
$ws.Cells.Item(35,3).Value = 214185
$ws.Cells.Item(35,4).Value = 93716
$ws.Cells.Item(35,5).Value = 93716
$ws.Cells.Item(35,6).Value = -913751
$ws.Cells.Item(35,7).Value = -90.7
$ws.Cells.Item(35,8).Value = -120469
$ws.Cells.Item(35,9).Value = -56.25

# Row 36: VA
$ws.Cells.Item(36,1).Value = "VA"
$ws.Cells.Item(36,2).Value = 403298
$ws.Cells.Item(36,3).Value = 265421
$ws.Cells.Item(36,4).Value = 1921763
$ws.Cells.Item(36,5).Value = 1426533
$ws.Cells.Item(36,6).Value = 1518465
$ws.Cells.Item(36,7).Value = 376.51
$ws.Cells.Item(36,8).Value = 1161112
$ws.Cells.Item(36,9).Value = 437.46

# Row 37: VT
$ws.Cells.Item(37,1).Value = "VT"
$ws.Cells.Item(37,2).Value = 62679
$ws.Cells.Item(37,3).Value = 40905
$ws.Cells.Item(37,4).Value = 424723
$ws.Cells.Item(37,5).Value = 151764
$ws.Cells.Item(37,6).Value = 362044
$ws.Cells.Item(37,7).Value = 577.62
$ws.Cells.Item(37,8).Value = 110859
$ws.Cells.Item(37,9).Value = 271.02

# Row 38: WI
$ws.Cells.Item(38,1).Value = "WI"
$ws.Cells.Item(38,2).Value = 478829
$ws.Cells.Item(38,3).Value = 419525
$ws.Cells.Item(38,4).Value = 1376055
$ws.Cells.Item(38,5).Value = 751288
$ws.Cells.Item(38,6).Value = 897226
$ws.Cells.Item(38,7).Value = 187.38
$ws.Cells.Item(38,8).Value = 331763
$ws.Cells.Item(38,9).Value = 79.08

# Row 39: WV
$ws.Cells.Item(39,1).Value = "WV"
$ws.Cells.Item(39,2).Value = 51412
$ws.Cells.Item(39,3).Value = 44368
$ws.Cells.Item(39,4).Value = 138943
$ws.Cells.Item(39,5).Value = 95038
$ws.Cells.Item(39,6).Value = 87531
$ws.Cells.Item(39,7).Value = 170.25
$ws.Cells.Item(39,8).Value = 50670
$ws.Cells.Item(39,9).Value = 114.2

# Row 40: WY
$ws.Cells.Item(40,1).Value = "WY"
$ws.Cells.Item(40,2).Value = 57610
$ws.Cells.Item(40,3).Value = 41737
$ws.Cells.Item(40,4).Value = 104263
$ws.Cells.Item(40,5).Value = 79037
$ws.Cells.Item(40,6).Value = 46653
$ws.Cells.Item(40,7).Value = 80.98
$ws.Cells.Item(40,8).Value = 37300
$ws.Cells.Item(40,9).Value = 89.37000000000001
